$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of original D,J,K,M,P values for rows 120..210 (index i -> old row 120+i)
$oldData = @(
    @(44505, 3200, 400, 450, 900),
    @(44168, 2600, 400, 450, 900),
    @(44179, 2400, 400, 450, 900),
    @(44510, 3000, 450, 475, 950),
    @(44161, 2600, 400, 450, 900),
    @(44477, 3300, 400, 450, 900),
    @(44438, 3460, 400, 450, 900),
    @(44498, 3200, 400, 450, 900),
    @(44413, 3140, 400, 450, 900),
    @(44160, 2500, 400, 450, 900),
    @(44308, 3000, 450, 475, 950),
    @(44431, 3400, 400, 450, 900),
    @(44357, 3000, 450, 475, 950),
    @(44218, 2800, 450, 475, 950),
    @(44335, 3100, 450, 475, 950),
    @(44251, 3000, 450, 475, 950),
    @(44295, 3300, 400, 450, 900),
    @(44210, 2800, 450, 475, 950),
    @(44407, 3500, 400, 450, 900),
    @(44343, 3000, 450, 475, 950),
    @(44230, 3000, 400, 450, 900),
    @(44316, 3360, 400, 450, 900),
    @(44265, 3200, 400, 450, 900),
    @(44329, 3060, 450, 475, 950),
    @(44186, 2400, 400, 450, 900),
    @(44526, 3320, 400, 450, 900),
    @(44272, 3100, 400, 450, 900),
    @(44167, 2400, 400, 450, 900),
    @(44211, 2600, 450, 475, 950),
    @(44489, 3000, 450, 475, 950),
    @(44209, 2600, 450, 475, 950),
    @(44389, 3120, 400, 450, 900),
    @(44517, 3000, 450, 475, 950),
    @(44405, 3100, 400, 450, 900),
    @(44221, 2900, 450, 475, 950),
    @(44280, 3000, 400, 450, 900),
    @(44330, 3340, 400, 450, 900),
    @(44483, 3060, 400, 450, 900),
    @(44448, 3000, 400, 450, 900),
    @(44196, 2800, 400, 450, 900),
    @(44463, 3400, 400, 450, 900),
    @(44239, 3000, 450, 475, 950),
    @(44476, 2800, 400, 450, 900),
    @(44169, 2900, 400, 450, 900),
    @(44496, 2900, 450, 475, 950),
    @(44328, 3040, 450, 475, 950),
    @(44515, 2200, 400, 450, 900),
    @(44560, 3080, 400, 450, 900),
    @(44543, 2300, 400, 450, 900),
    @(44552, 3000, 450, 475, 950),
    @(44379, 3600, 400, 450, 900),
    @(44216, 2600, 450, 475, 950),
    @(44426, 3200, 450, 475, 950),
    @(44286, 3000, 400, 450, 900),
    @(44452, 3200, 400, 450, 900),
    @(44473, 2000, 400, 450, 900),
    @(44298, 2800, 400, 450, 900),
    @(44482, 2960, 450, 475, 950),
    @(44294, 3000, 400, 450, 900),
    @(44554, 2800, 400, 450, 900),
    @(44305, 2800, 400, 450, 900),
    @(44545, 2800, 450, 475, 950),
    @(44445, 3400, 400, 450, 900),
    @(44301, 3000, 400, 450, 900),
    @(44193, 3000, 400, 450, 900),
    @(44454, 3200, 450, 475, 950),
    @(44466, 3320, 400, 450, 900),
    @(44392, 3000, 400, 450, 900),
    @(44412, 3200, 400, 450, 900),
    @(44322, 3080, 450, 475, 950),
    @(44162, 2800, 400, 450, 900),
    @(44434, 3140, 400, 450, 900),
    @(44532, 3000, 400, 450, 900),
    @(44490, 3000, 400, 450, 900),
    @(44427, 3140, 400, 450, 900),
    @(44491, 3400, 400, 450, 900),
    @(44293, 3100, 400, 450, 900),
    @(44266, 3200, 450, 475, 950),
    @(44533, 3360, 400, 450, 900),
    @(44494, 2200, 400, 450, 900),
    @(44571, 2400, 400, 450, 900),
    @(44279, 3200, 400, 450, 900),
    @(44277, 2800, 400, 450, 900),
    @(44525, 2000, 400, 450, 900),
    @(44354, 3000, 400, 450, 900),
    @(44503, 2700, 450, 475, 950),
    @(44462, 3000, 400, 450, 900),
    @(44384, 3100, 450, 475, 950),
    @(44512, 3340, 400, 450, 900),
    @(44312, 3000, 400, 450, 900),
    @(44511, 3100, 400, 450, 900)
)

# Row 120 becomes a brand new data point (date 44574, volume 2900).
# Its Precio minimo / Precio promedio ponderado / Precio $/Kg (K,M,P) stay
# the same as before the edit (400, 450, 900), so only D and J change.
$ws.Cells.Item(120, 4).Value = 44574
$ws.Cells.Item(120, 10).Value = 2900

# Rows 121..210 each take on the D/J/K/M/P values that used to belong to
# the row immediately above them (a "shift down by one" over the old
# row-120..210 block), pushing the old row-210 values into a brand new
# row 211.
for ($i = 1; $i -le 90; $i++) {
    $destRow = 120 + $i
    $src = $oldData[$i - 1]
    $ws.Cells.Item($destRow, 4).Value  = $src[0]   # D - Fecha
    $ws.Cells.Item($destRow, 10).Value = $src[1]   # J - Volumen
    $ws.Cells.Item($destRow, 11).Value = $src[2]   # K - Precio minimo
    $ws.Cells.Item($destRow, 13).Value = $src[3]   # M - Precio promedio ponderado
    $ws.Cells.Item($destRow, 16).Value = $src[4]   # P - Precio $/Kg
}

# Row 211 is entirely new, so besides D/J/K/M/P (taken from the old row
# 210) it also needs the constant "template" columns copied across, the
# same way every other row in this data block is populated.
$ws.Cells.Item(211, 1).Value  = 8
$ws.Cells.Item(211, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(211, 3).Value  = "Coquimbo"
$ws.Cells.Item(211, 5).Value  = 4
$ws.Cells.Item(211, 6).Value  = 100112012
$ws.Cells.Item(211, 7).Value  = "Espinaca"
$ws.Cells.Item(211, 8).Value  = "Sin especificar"
$ws.Cells.Item(211, 9).Value  = "Primera"
$ws.Cells.Item(211, 12).Value = 500
$ws.Cells.Item(211, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(211, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(211, 17).Value = 0.5
$ws.Cells.Item(211, 18).Value = "Hortaliza"

$last = $oldData[90]
$ws.Cells.Item(211, 4).Value  = $last[0]
$ws.Cells.Item(211, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(211, 10).Value = $last[1]
$ws.Cells.Item(211, 11).Value = $last[2]
$ws.Cells.Item(211, 13).Value = $last[3]
$ws.Cells.Item(211, 16).Value = $last[4]
